$wb = $excel.ActiveWorkbook

# --- Model identifier renames (sheet tab names track the adjusted model sizes) ---
$wb.Worksheets.Item(1).Name = "iCC389"
$wb.Worksheets.Item(3).Name = "iCC470"
$wb.Worksheets.Item(4).Name = "iCC651"

# --- Updated totals for the adjusted iCC651 model ---
$ws4 = $wb.Worksheets.Item("iCC651")
$ws4.Range("B2").Value = 269
$ws4.Range("C2").Value = 740
$ws4.Range("D2").Value = 136
$ws4.Range("F2").Value = 137
$ws4.Range("H2").Value = 681

# Re-align iCC651's header/label formatting with the shared style used by the
# other sheets (the old per-sheet duplicate style is dropped).
$ws1 = $wb.Worksheets.Item("iCC389")
$ws1.Range("B1:H1").Copy()
$ws4.Range("B1:H1").PasteSpecial(-4122)  # xlPasteFormats
$ws1.Range("A2").Copy()
$ws4.Range("A2").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Reset stale selections back to A1 and make the first sheet active ---
$ws3 = $wb.Worksheets.Item("iCC470")
$ws2 = $wb.Worksheets.Item(2)

[void]$ws2.Range("A1").Select()
[void]$ws3.Range("A1").Select()
[void]$ws4.Range("A1").Select()
[void]$ws1.Range("A1").Select()
